# x1469 Humans have been renamed "Homo sapiens (Human)"
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet3")

# Every cell that previously read "Human" now reads "Homo sapiens (Human)"
$ws.Range("F3").Value = "Homo sapiens (Human)"
$ws.Range("F5").Value = "Homo sapiens (Human)"

# Selection moved from J5 to F5
$ws.Activate()
$ws.Range("F5").Select()
